$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.087.36"
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = "'2.371.66"
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'303.08"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").Value = "'95.39"
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'0.502"
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").Value = "'0.481"
$ws.Range("E9").Value = '  -2.81%  '
$ws.Range("D10").Value = "'34.30"
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = "'0.125"
$ws.Range("E11").Value = '  +2.62%  '
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = "'18.33"
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("D14").Value = "'6.75"
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = "'2.738.50"
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = "'2.384.26"
$ws.Range("E16").Value = '  +2.13%  '
$ws.Range("D17").Value = "'0.798"
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = "'43.103.05"
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").Value = "'11.97"
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = "'235.62"
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("E28").Value = '  +15.23%  '
$ws.Range("D29").Value = "'9.35"
$ws.Range("E29").Value = '  +2.46%  '
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").Value = "'17.61"
$ws.Range("E33").Value = '  +2.35%  '
$ws.Range("E34").Value = '  +9.27%  '
$ws.Range("D35").Value = "'0.0730"
$ws.Range("E35").Value = '  -3.96%  '
$ws.Range("E36").Value = '  +1.22%  '
$ws.Range("D37").Value = "'127.25"
$ws.Range("E37").Value = '  +1.31%  '
$ws.Range("E38").Value = '  +4.37%  '
$ws.Range("E39").Value = '  -1.93%  '
$ws.Range("E40").Value = '  -2.45%  '
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").Value = "'20.74"
$ws.Range("E42").Value = '  -5.20%  '
$ws.Range("D43").Value = "'1.929.88"
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("D46").Value = "'2.74"
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("D47").Value = "'9.19"
$ws.Range("E47").Value = '  -8.85%  '
$ws.Range("D48").Value = "'2.598.65"
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("E49").Value = '  +2.65%  '
$ws.Range("D50").Value = "'71.44"
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("E51").Value = '  +1.38%  '
